$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.450.20'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.80%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.777.55'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.66%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.67'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.41'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.38%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.776.54'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.67%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.28%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.58%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.62%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000266'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.67'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.95%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.419.55'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.83%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.766.26'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.64%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.522.46'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.96%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.40'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.44%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.06'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.81%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.98'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '468.55'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.74%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.703'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.99%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.21'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000143'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.24'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.40%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.16'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.930.36'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.81'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.41'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.87%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.23'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.16%  '

$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.13'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.34'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.08%  '


$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.737.49'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.51%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.97%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.46'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -7.81%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.139'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.01'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.26%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.19%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.82%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.97'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.20%  '

$ws.Range("B47").Value = 'Cosmos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.59'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.54%  '

$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '43.10'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +12.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '403.69'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.78'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.72%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '145.88'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.08%  '
